# Added aggregate row (AggregateMapping) to mapping spreadsheet and pie
# chart age analysis.
#
# Adds a 4th column (D) to Sheet1 that buckets each bggrecagerange row
# into one of four coarse life-stage groups, derived from the existing
# numeric SortMapping value in column C:
#   C <= 7   -> "Kids"
#   C 8-12   -> "Preteens"
#   C 13-17  -> "Teens"
#   C >= 18  -> "Adults"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header -----------------------------------------------------------
$ws.Range("D1").Value = "AggregateMapping"
$ws.Range("D1").Font.Bold = $true

# --- Body: bucket each data row (2-62) by its column-C numeric value --
$kids     = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32)
$preteens = @(33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48)
$teens    = @(49,50,51,52,53,54,55,56)
$adults   = @(57,58,59,60,61,62)

foreach ($r in $kids)     { $ws.Cells.Item($r, 4).Value = "Kids" }
foreach ($r in $preteens) { $ws.Cells.Item($r, 4).Value = "Preteens" }
foreach ($r in $teens)    { $ws.Cells.Item($r, 4).Value = "Teens" }
foreach ($r in $adults)   { $ws.Cells.Item($r, 4).Value = "Adults" }

# --- Column widths for the new columns ---------------------------------
# Stored widths in the xlsx are target + 5/6 (Excel's column-width padding),
# so subtract that back off to land exactly on the target stored widths.
$ws.Range("C1").ColumnWidth = 13.5 - 0.8333333333333334
$ws.Range("D1").ColumnWidth = 20 - 0.8333333333333334

# --- View: zoom in and select the newly added Adults block -------------
$excel.ActiveWindow.Zoom = 160
$null = $ws.Range("D57:D62").Select()
